# Update the "Förändrad" (Changed) date column (C) from 2023-09-11 to
# 2023-09-12 for every data row (rows 2 through 499) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 499
$rng = $ws.Range("C2:C" + $lastRow)

foreach ($cell in $rng.Cells) {
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
